$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 15: "6 - 12 toukokuuta" (introduces shared string BEFORE the
# "26.3? + -" edit below, so the shared-string table ends up in the same
# order as the target file: 12=toukokuuta, 13=26.3?+-, 14=bootstrappia) ---
$ws.Range("A15").NumberFormat = "mm-dd-yy"
$ws.Range("A15").Value = "6 - 12 toukokuuta"

# Row 2: update the note text and hours worked
$ws.Range("A2").Value = "26.3? + -"
$ws.Range("B2").Value = 10

# Finish new row 15
$ws.Range("B15").Value = 10
$ws.Range("C15").Value = "koodin ja MHn sivujen tutkimista kotona + pikkusen ehkä bootstrappia"

# New row 16: a dated entry
$ws.Range("A16").NumberFormat = "d-mmm"
$ws.Range("A16").Value = 41407
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = "projektihuoneella"

# Update the saved selection to match where the user ended up
$ws.Range("C17").Select()
